$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C68").Value = "11-23755885"
$ws.Range("C65").Value = "11-76381379"
$ws.Range("C64").Value = "11-63620357"

$ws.Range("C65").Select()
